$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row "A2.3.1 Ja " (old row 17) was removed entirely; deleting the row
# shifts every row below it up by one and keeps their existing formatting.
$ws.Rows(17).Delete()

# --- Fix up ID (col A) and ParentID (col B) values that were left stale
# after the shift, and apply the remaining content edits from the diff ---

# Row 15: A2.2 Flaechenaufstellung -> ParentID now 11 (A2. Betriebsstruktur),
# and gains an "Upload" = Ja value (new D cell, needs text formatting).
$ws.Cells.Item(15, 2).Value = "11"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "Ja"

# Row 16: A2.3 Tierhaltung -> gains an "Upload" = ja value (new D cell).
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "ja"

# Row 17: A2.4 Waldflaechen -> ID 16, ParentID 11
$ws.Cells.Item(17, 1).Value = "16"
$ws.Cells.Item(17, 2).Value = "11"

# Row 18: merged question text, ID 17 (ParentID stays 17 - unchanged)
$ws.Cells.Item(18, 1).Value = "17"
$ws.Cells.Item(18, 3).Value = "A2.4 Waldflächen: Hektar & typische Nutzungen angeben "

# Row 19: A3. Nebetaetigkeiten gesamt -> ID 18, ParentID 17
$ws.Cells.Item(19, 1).Value = "18"
$ws.Cells.Item(19, 2).Value = "17"

# Rows 20-25: A3.1 .. A3.6 -> ID shifts down by one, ParentID becomes 18
$ws.Cells.Item(20, 1).Value = "19"
$ws.Cells.Item(20, 2).Value = "18"

$ws.Cells.Item(21, 1).Value = "20"
$ws.Cells.Item(21, 2).Value = "18"

$ws.Cells.Item(22, 1).Value = "21"
$ws.Cells.Item(22, 2).Value = "18"

$ws.Cells.Item(23, 1).Value = "22"
$ws.Cells.Item(23, 2).Value = "18"

$ws.Cells.Item(24, 1).Value = "23"
$ws.Cells.Item(24, 2).Value = "18"

$ws.Cells.Item(25, 1).Value = "24"
$ws.Cells.Item(25, 2).Value = "18"

# Row 26: B1. Stammdaten & Betreiber -> ID 25, ParentID cleared entirely
$ws.Cells.Item(26, 1).Value = "25"
$ws.Cells.Item(26, 2).Clear()

# Rows 27-33: only the ID (col A) shifts down by one
$ws.Cells.Item(27, 1).Value = "26"
$ws.Cells.Item(28, 1).Value = "27"
$ws.Cells.Item(29, 1).Value = "28"
$ws.Cells.Item(30, 1).Value = "29"
$ws.Cells.Item(31, 1).Value = "30"
$ws.Cells.Item(32, 1).Value = "31"
$ws.Cells.Item(33, 1).Value = "32"

# Restore the originally-selected cell (it shifted up together with the
# deleted row, from B27 to B26).
[void]$ws.Range("B26").Select()
